$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "пока" / "давай, до свидания"  ->  "до свидания" / "и вам не хворать"
$ws.Range("A3").Value = "до свидания"
$ws.Range("C3").Value = "и вам не хворать"

# Row 4 ("плач" / file_id) is removed entirely
$ws.Rows.Item(4).Delete()

# Move the active selection to A6 (matches the saved view state)
$null = $ws.Range("A6").Select()
